$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "Observed" (F1) and "Expected" (G1) columns, shifting the
# existing "Is Significant" header into the new H1 cell. ---
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"
$ws.Range("H1").Value = "Is Significant"

# Give the two brand-new header cells the same look (bold, bordered, centered)
# as the rest of the header row by copying the formatting from F1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# PasteSpecial(formats) does not touch cell contents, but make sure the text is
# still correct afterwards.
$ws.Range("G1").Value = "Expected"
$ws.Range("H1").Value = "Is Significant"

# --- Row 2 / Row 3 data: observed & expected contingency arrays, "Is
# Significant" flag moves out to the new H column. ---
$ws.Range("F2").Value = "[240 358] ; [16  4]"
$ws.Range("F3").Value = "[269 329] ; [17  3]"

$ws.Range("G2").Value = "[247.71521036 350.28478964] ; [ 8.28478964 11.71521036]"
$ws.Range("G3").Value = "[276.74433657 321.25566343] ; [ 9.25566343 10.74433657]"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
